$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 44754
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = 100112043
$ws.Range("G10").Value = "Pepino dulce"
$ws.Range("H10").Value = "Cultivar IV Región"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 80
$ws.Range("K10").Value = 16000
$ws.Range("L10").Value = 17000
$ws.Range("M10").Value = 16500
$ws.Range("N10").Value = "`$/bandeja 18 kilos"
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 917
$ws.Range("Q10").Value = 18
$ws.Range("R10").Value = "Hortaliza"
